# Add a new worksheet "4D-merges" right after the existing "5D" sheet.
# It starts life as a copy of "5D" (same column widths / row heights /
# per-cell styling) and we then rewrite the handful of cells that differ,
# matching the new example's template strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("5D")

# Duplicate the "5D" sheet, placing the copy immediately after it.
$ws.Copy($null, $ws)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "4D-merges"

# Row 2: A2 keeps a header-style placeholder, but now clarifies it targets
# the second element of the "data" array. B2 takes over the template that
# used to live in C2, with an explicit "!mergeCells=false" option. C2 is
# cleared since the value moved to B2.
$ws2.Range("A2").Value = "{{ | data.1.data | header | 1:0 }}"
$ws2.Range("B2").Value = "{{ A2 | 1 * data | header | 1 | !mergeCells=false}}"
$ws2.Range("C2").ClearContents()

# Row 3: B3 takes over a (rewritten) template that used to live in B3/C2.
# C3 is cleared since its old content isn't used on this sheet.
$ws2.Range("B3").Value = "{{ B2 | data * data || 1 }}"
$ws2.Range("C3").ClearContents()
